$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns of affected rows stay text-typed (matches source inlineStr cells),
# so numeric-looking values like "1.008" or "12.70" are not coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = "30.496.08"
$ws.Range("E2").Value2 = "  +0.23%  "
$ws.Range("D3").Value2 = "2.134.36"
$ws.Range("E3").Value2 = "  +1.60%  "
$ws.Range("D4").Value2 = "1.008"
$ws.Range("E4").Value2 = "  +0.52%  "
$ws.Range("D5").Value2 = "352.32"
$ws.Range("E5").Value2 = "  +5.27%  "
$ws.Range("E6").Value2 = "  +0.45%  "
$ws.Range("D7").Value2 = "0.5253"
$ws.Range("E7").Value2 = "  +0.69%  "
$ws.Range("D8").Value2 = "0.4545"
$ws.Range("E8").Value2 = "  -0.15%  "
$ws.Range("D9").Value2 = "53.71"
$ws.Range("E9").Value2 = "  -0.89%  "
$ws.Range("D10").Value2 = "0.09130"
$ws.Range("E10").Value2 = "  +2.71%  "
$ws.Range("D11").Value2 = "1.192"
$ws.Range("E11").Value2 = "  +1.18%  "
$ws.Range("D12").Value2 = "25.52"
$ws.Range("E12").Value2 = "  +5.69%  "
$ws.Range("D13").Value2 = "2.139.29"
$ws.Range("E13").Value2 = "  +2.35%  "
$ws.Range("D14").Value2 = "6.877"
$ws.Range("E14").Value2 = "  +1.23%  "
$ws.Range("D15").Value2 = "8.151"
$ws.Range("E15").Value2 = "  +1.73%  "
$ws.Range("D16").Value2 = "101.50"
$ws.Range("E16").Value2 = "  +4.59%  "
$ws.Range("D17").Value2 = "0.00001166"
$ws.Range("E17").Value2 = "  +1.81%  "
$ws.Range("E18").Value2 = "  +0.41%  "
$ws.Range("B19").Value2 = "Avalanche"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value2 = "20.61"
$ws.Range("E19").Value2 = "  +7.50%  "
$ws.Range("B20").Value2 = "TRON"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value2 = "0.06711"
$ws.Range("E20").Value2 = "  +1.35%  "
$ws.Range("E21").Value2 = "  +0.41%  "
$ws.Range("D22").Value2 = "6.367"
$ws.Range("E22").Value2 = "  +1.17%  "
$ws.Range("D23").Value2 = "30.597.62"
$ws.Range("E23").Value2 = "  +0.39%  "
$ws.Range("D24").Value2 = "12.86"
$ws.Range("E24").Value2 = "  +4.27%  "
$ws.Range("D25").Value2 = "2.384"
$ws.Range("E25").Value2 = "  +1.15%  "
$ws.Range("D26").Value2 = "2.393.65"
$ws.Range("E26").Value2 = "  +2.45%  "
$ws.Range("E27").Value2 = "  +1.39%  "
$ws.Range("D28").Value2 = "2.612"
$ws.Range("E28").Value2 = "  +3.98%  "
$ws.Range("D29").Value2 = "164.96"
$ws.Range("E29").Value2 = "  +1.60%  "
$ws.Range("D30").Value2 = "135.74"
$ws.Range("E30").Value2 = "  +1.97%  "
$ws.Range("D31").Value2 = "1.217"
$ws.Range("E31").Value2 = "  +1.08%  "
$ws.Range("B32").Value2 = "Stellar"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value2 = "0.1082"
$ws.Range("E32").Value2 = "  +1.44%  "
$ws.Range("B33").Value2 = "ARBITRUM"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value2 = "1.709"
$ws.Range("E33").Value2 = "  +3.47%  "
$ws.Range("D34").Value2 = "6.384"
$ws.Range("E34").Value2 = "  -0.19%  "
$ws.Range("D35").Value2 = "4.031"
$ws.Range("E35").Value2 = "  +2.51%  "
$ws.Range("D36").Value2 = "6.122"
$ws.Range("E36").Value2 = "  +4.68%  "
$ws.Range("E37").Value2 = "  +1.13%  "
$ws.Range("E38").Value2 = "  +2.85%  "
$ws.Range("D39").Value2 = "0.06971"
$ws.Range("E39").Value2 = "  +1.90%  "
$ws.Range("D40").Value2 = "0.2353"
$ws.Range("E40").Value2 = "  +1.70%  "
$ws.Range("D41").Value2 = "12.70"
$ws.Range("E41").Value2 = "  +0.36%  "
$ws.Range("D42").Value2 = "0.6979"
$ws.Range("E42").Value2 = "  +1.62%  "
$ws.Range("D43").Value2 = "1.273"
$ws.Range("E43").Value2 = "  +1.87%  "
$ws.Range("D44").Value2 = "14.79"
$ws.Range("E44").Value2 = "  +5.75%  "
$ws.Range("D45").Value2 = "0.6509"
$ws.Range("E45").Value2 = "  +1.85%  "
$ws.Range("D46").Value2 = "2.348"
$ws.Range("E46").Value2 = "  +1.22%  "
$ws.Range("D47").Value2 = "0.00000000376"
$ws.Range("E47").Value2 = "  +11.29%  "
$ws.Range("D48").Value2 = "3.744"
$ws.Range("E48").Value2 = "  +2.23%  "
$ws.Range("D49").Value2 = "1.248"
$ws.Range("E49").Value2 = "  +0.12%  "
$ws.Range("D50").Value2 = "83.69"
$ws.Range("E50").Value2 = "  +0.75%  "
$ws.Range("E51").Value2 = "  +2.28%  "

# Restore default (unstyled) cell style now that values are set, matching the original workbook
# which had no explicit style index on these data cells.
$ws.Range("D2:E51").Style = "Normal"
